$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.124902828213514
$ws.Range("K2").Value = 0.0642794618386348
$ws.Range("L2").Value = 0.118423810913776
$ws.Range("N2").Value = 0.0475638500827432
